# Applies the edits described by the diff:
#  - Metadata sheet: Date and FHIR Version values updated
#  - Elements sheet: ele-1 constraint text simplified (Extension row),
#    Extension.id Type(s) changed from "id" to "string",
#    Extension.value[x] Type(s) list updated (drop CodeableReference, de-dup Ratio/RatioRange, add Meta),
#    Extension.value[x] Definition text updated from R4B to R4 link

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$wsMeta.Range("B15").Value = "4.0.1"

$wsElem = $wb.Worksheets.Item("Elements")

$wsElem.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

$wsElem.Range("K3").Value = "string`n"

$wsElem.Range("K8").Value = "base64Binary`nbooleancanonicalcodedatedateTimedecimalidinstantintegermarkdownoidpositiveIntstringtimeunsignedInturiurluuidAddressAgeAnnotationAttachmentCodeableConceptCodingContactPointCountDistanceDurationHumanNameIdentifierMoneyPeriodQuantityRangeRatioReferenceSampledDataSignatureTimingContactDetailContributorDataRequirementExpressionParameterDefinitionRelatedArtifactTriggerDefinitionUsageContextDosageMeta"

$wsElem.Range("M8").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
